$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1. Drop the old "_GoBack" bookmark (it sat between "execução de" and
#    " determinadas tarefas" in the original doc) -- it is replaced by the
#    new named bookmarks added below.
# ---------------------------------------------------------------------------
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks("_GoBack").Delete()
}

# ---------------------------------------------------------------------------
# 2. Fix the missing word: "devido perda de tempo" -> "devido a perda de tempo"
# ---------------------------------------------------------------------------
$d.Content.Find.Execute("devido perda de tempo", $true, $false, $false, $false, `
    $false, $true, 1, $false, "devido a perda de tempo", 2)

# ---------------------------------------------------------------------------
# 3. Wrap specific phrases with the new "_Hlk…" bookmarks (left to right, so
#    the w:id values come out 0,1,2,3 in document order).
# ---------------------------------------------------------------------------
$r = $d.Content
$r.Find.Execute("dificuldade de realizar tarefas comuns dentro das residências")
$d.Bookmarks.Add("_Hlk41036457", $r)

$r = $d.Content
$r.Find.Execute("falta de praticidade para execução de determinadas tarefas")
$d.Bookmarks.Add("_Hlk41037245", $r)

$r = $d.Content
$r.Find.Execute("e falta de comodidade")
$d.Bookmarks.Add("_Hlk41037403", $r)

$r = $d.Content
$r.Find.Execute("pessoas indispostas")
$paraEnd = $r.Paragraphs(1).Range.End
$r.SetRange($r.Start, $paraEnd)
$d.Bookmarks.Add("_Hlk41037151", $r)

Write-Output "edits applied"
